# feat: add 2022-Q1 data
#
# Target layout:
#   2021-Q4  (sheetId 1, unchanged)
#   2022-Q1  (sheetId 2, NEW - fund holding table, same border/bold style
#             that "总计" used)
#   总计      (sheetId 3, same sheet as before but with a new first data
#             row for 2022-Q1; old 2021-Q4 total row pushed to row 3)
#
# Approach: duplicate "总计" (so the new sheet inherits its style index)
# and place the duplicate right after the original - then swap names so
# the ORIGINAL keeps sheetId 2 (becomes "2022-Q1") and the COPY takes
# sheetId 3 (becomes "总计"). That reproduces the sheetId numbering in
# the target exactly.

$wb = $excel.ActiveWorkbook
$q4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")

$total.Copy($null, $total)
$totalCopy = $wb.Worksheets.Item(3)

$q1 = $total
$q1.Name = "2022-Q1"
$totalCopy.Name = "总计"

# --- "2022-Q1": rebuild as the fund-holding table (columns A-H) ---

# New header cells E1:H1 inherit the "总计" header style (s=2) by
# copying format from the existing D1 header cell first.
$q1.Range("D1").Copy($q1.Range("E1:H1"))

$q1.Range("B1").Value2 = "基金代码"
$q1.Range("C1").Value2 = "基金名称"
$q1.Range("D1").Value2 = "基金规模"
$q1.Range("E1").Value2 = "股票总仓位"
$q1.Range("F1").Value2 = "仓位占比"
$q1.Range("G1").Value2 = "持有市值(亿元)"
$q1.Range("H1").Value2 = "仓位排名"

# Row 2 - 011444 创金合信瑞裕混合A
$q1.Range("A2").Value2 = 0
$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value2 = "011444"
$q1.Range("C2").Value2 = "创金合信瑞裕混合A"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value2 = "0.03"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value2 = "68.73"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value2 = "2.87"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value2 = "0.0009"
$q1.Range("H2").Value2 = 7

# Row 3 - 011445 创金合信瑞裕混合C
# A3 is a brand-new cell (the "总计" template only had 2 rows), so copy
# A2's style (s=2) onto it before writing the value.
$q1.Range("A2").Copy($q1.Range("A3"))
$q1.Range("A3").Value2 = 1
$q1.Range("B3").NumberFormat = "@"
$q1.Range("B3").Value2 = "011445"
$q1.Range("C3").Value2 = "创金合信瑞裕混合C"
$q1.Range("D3").NumberFormat = "@"
$q1.Range("D3").Value2 = "0.00"
$q1.Range("E3").NumberFormat = "@"
$q1.Range("E3").Value2 = "68.73"
$q1.Range("F3").NumberFormat = "@"
$q1.Range("F3").Value2 = "2.87"
$q1.Range("G3").Value2 = 0
$q1.Range("H3").Value2 = 7

# --- "总计": shift the existing total row down, add the 2022-Q1 row ---

$totalCopy.Range("A2").Copy($totalCopy.Range("A3"))
$totalCopy.Range("A3").Value2 = 1
$totalCopy.Range("B3").Value2 = "2021-Q4"
$totalCopy.Range("C3").Value2 = 2
$totalCopy.Range("D3").Value2 = 0

$totalCopy.Range("A2").Value2 = 0
$totalCopy.Range("B2").Value2 = "2022-Q1"
$totalCopy.Range("C2").Value2 = 2
$totalCopy.Range("D2").Value2 = 0

# Restore the original active sheet / selection (copying activates the
# new sheet, which would otherwise steal "tabSelected").
$q4.Activate()
$q4.Range("A1").Select() | Out-Null
